{"js": "// Remove the blank paragraph, the \"Ver no Jupiter...\" paragraph, and the\n// \"\u00a9 2020 ...\" copyright paragraph that used to sit right after the\n// \"LOM3057: ...\" requirement line, leaving the trailing blank paragraph\n// (the one right before the page-break paragraph) untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst markerText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n\nconst removeIdx = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === markerText) {\n    // The blank paragraph right before this one, this one, and the\n    // copyright paragraph right after it are all removed (3 paragraphs).\n    removeIdx.push(i - 1, i, i + 1);\n  }\n}\n\n// Delete from the highest index down so earlier indices stay valid.\nremoveIdx\n  .sort((a, b) => b - a)\n  .forEach((idx) => {\n    if (idx >= 0 && idx < paragraphs.items.length) {\n      paragraphs.items[idx].delete();\n    }\n  });\n\nawait context.sync();\n", "ps1": "# Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n# \"\u00a9 2020 ...\" paragraph that used to follow the LOM3057 requirement line,\n# leaving the trailing blank paragraph (right before the page-break\n# paragraph) untouched.\n$d = $word.ActiveDocument\n\n$markerText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $markerText) {\n        # Delete this paragraph plus the copyright paragraph right after it\n        # and the blank paragraph right before it (3 paragraphs total).\n        $d.Paragraphs.Item($i + 1).Range.Delete()\n        $d.Paragraphs.Item($i).Range.Delete()\n        $d.Paragraphs.Item($i - 1).Range.Delete()\n        break\n    }\n}\n"}
